# Auto-generated PowerShell (Excel COM-interop) script
# Applies the diff: splits the single 'ODI Batting' sheet into three sheets:
#   1) Player Info          - new sheet, player bio fields
#   2) ODI Batting           - existing sheet, MATCH_CARD_LINK -> MATCH_CODE column rework
#   3) ODI Batting Extra     - new sheet, extra per-match batting stats

$wb = $excel.ActiveWorkbook

# --- Reference the original (pre-existing) sheet by name; it holds the ODI Batting data ---
$origName = $wb.Worksheets.Item(1).Name
$batting = $wb.Worksheets.Item($origName)

# --- Insert the new 'Player Info' sheet BEFORE the batting sheet ---
$playerInfo = $wb.Worksheets.Add($batting)
$playerInfo.Name = "Player Info"

# --- Insert the new 'ODI Batting Extra' sheet AFTER the batting sheet ---
$battingExtra = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($origName))
$battingExtra.Name = "ODI Batting Extra"

# --- Finally rename the original sheet to its target name ---
$wb.Worksheets.Item($origName).Name = "ODI Batting"
$batting = $wb.Worksheets.Item("ODI Batting")

## =============== Player Info ===============
$hdr = $playerInfo.Range("A1:D1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1

$arr = New-Object 'object[,]' 2,4
$arr[0,0] = 'ID'
$arr[0,1] = 'NAME'
$arr[0,2] = 'BATTING_HAND'
$arr[0,3] = 'BOWL_STYLE'
$arr[1,0] = '4412'
$arr[1,1] = 'Mohammad Rizwan'
$arr[1,2] = 'Right Handed'
$arr[1,3] = 'Does Not Bowl | Unknown'
$playerInfo.Range("A1:D2").NumberFormat = "@"
$playerInfo.Range("A1:D2").Value = $arr

## =============== ODI Batting (in-place edits) ===============
$batting.Range("D1").Value = "MATCH_CODE"

# MATCH_CODE data column (D2:D53) must stay textual, not auto-converted to numbers
$batting.Range("D2:D53").NumberFormat = "@"
$codeArr = New-Object 'object[,]' 52,1
$codeArr[0,0] = '3797'
$codeArr[1,0] = '3798'
$codeArr[2,0] = '3799'
$codeArr[3,0] = '3814'
$codeArr[4,0] = '3819'
$codeArr[5,0] = '3820'
$codeArr[6,0] = '3821'
$codeArr[7,0] = '3822'
$codeArr[8,0] = '3836'
$codeArr[9,0] = '3837'
$codeArr[10,0] = '3838'
$codeArr[11,0] = '3859'
$codeArr[12,0] = '3861'
$codeArr[13,0] = '3863'
$codeArr[14,0] = '3883'
$codeArr[15,0] = '3930'
$codeArr[16,0] = '3932'
$codeArr[17,0] = '3939'
$codeArr[18,0] = '3943'
$codeArr[19,0] = '3944'
$codeArr[20,0] = '3972'
$codeArr[21,0] = '3973'
$codeArr[22,0] = '3975'
$codeArr[23,0] = '3977'
$codeArr[24,0] = '3981'
$codeArr[25,0] = '4244'
$codeArr[26,0] = '4247'
$codeArr[27,0] = '4273'
$codeArr[28,0] = '4274'
$codeArr[29,0] = '4275'
$codeArr[30,0] = '4276'
$codeArr[31,0] = '4277'
$codeArr[32,0] = '4432'
$codeArr[33,0] = '4433'
$codeArr[34,0] = '4434'
$codeArr[35,0] = '4458'
$codeArr[36,0] = '4459'
$codeArr[37,0] = '4460'
$codeArr[38,0] = '4472'
$codeArr[39,0] = '4473'
$codeArr[40,0] = '4476'
$codeArr[41,0] = '4564'
$codeArr[42,0] = '4565'
$codeArr[43,0] = '4567'
$codeArr[44,0] = '4586'
$codeArr[45,0] = '4590'
$codeArr[46,0] = '4592'
$codeArr[47,0] = '4634'
$codeArr[48,0] = '4638'
$codeArr[49,0] = '4686'
$codeArr[50,0] = '4688'
$codeArr[51,0] = '4690'
$batting.Range("D2:D53").Value = $codeArr

# Clear INNING_NUMBER (column B) for 'did not bat' rows so the cell becomes blank/absent
$batting.Range("B8").Value = ""
$batting.Range("B12").Value = ""
$batting.Range("B23").Value = ""
$batting.Range("B29").Value = ""
$batting.Range("B45").Value = ""

## =============== ODI Batting Extra ===============
$hdr3 = $battingExtra.Range("A1:F1")
$hdr3.Font.Bold = $true
$hdr3.HorizontalAlignment = -4108
$hdr3.VerticalAlignment = -4160
$hdr3.Borders.LineStyle = 1

$arr3 = New-Object 'object[,]' 21,6
$arr3[0,0] = 'MATCH_CODE'
$arr3[0,1] = 'BATTING_POSITION'
$arr3[0,2] = 'NUM_4'
$arr3[0,3] = 'NUM_6'
$arr3[0,4] = 'PERCENT_RUNS_OF_TOTAL'
$arr3[0,5] = 'MAN_OF_MATCH'
$arr3[1,0] = '4432'
$arr3[1,1] = ""
$arr3[1,2] = ''
$arr3[1,3] = ''
$arr3[1,4] = ''
$arr3[1,5] = 'NO'
$arr3[2,0] = '4433'
$arr3[2,1] = ""
$arr3[2,2] = ''
$arr3[2,3] = ''
$arr3[2,4] = ''
$arr3[2,5] = 'NO'
$arr3[3,0] = '4434'
$arr3[3,1] = 5
$arr3[3,2] = '0'
$arr3[3,3] = '0'
$arr3[3,4] = '3.60%'
$arr3[3,5] = 'NO'
$arr3[4,0] = '4458'
$arr3[4,1] = 4
$arr3[4,2] = '4'
$arr3[4,3] = '0'
$arr3[4,4] = '14.60%'
$arr3[4,5] = 'NO'
$arr3[5,0] = '4459'
$arr3[5,1] = 4
$arr3[5,2] = '0'
$arr3[5,3] = '0'
$arr3[5,4] = ''
$arr3[5,5] = 'NO'
$arr3[6,0] = '4460'
$arr3[6,1] = 4
$arr3[6,2] = '0'
$arr3[6,3] = '0'
$arr3[6,4] = '0.63%'
$arr3[6,5] = 'NO'
$arr3[7,0] = '4472'
$arr3[7,1] = 4
$arr3[7,2] = '3'
$arr3[7,3] = '0'
$arr3[7,4] = '9.22%'
$arr3[7,5] = 'NO'
$arr3[8,0] = '4473'
$arr3[8,1] = 4
$arr3[8,2] = '1'
$arr3[8,3] = '0'
$arr3[8,4] = '2.56%'
$arr3[8,5] = 'NO'
$arr3[9,0] = '4476'
$arr3[9,1] = 4
$arr3[9,2] = '8'
$arr3[9,3] = '0'
$arr3[9,4] = '22.36%'
$arr3[9,5] = 'NO'
$arr3[10,0] = '4564'
$arr3[10,1] = 5
$arr3[10,2] = '0'
$arr3[10,3] = '0'
$arr3[10,4] = '4.44%'
$arr3[10,5] = 'NO'
$arr3[11,0] = '4565'
$arr3[11,1] = ""
$arr3[11,2] = ''
$arr3[11,3] = ''
$arr3[11,4] = ''
$arr3[11,5] = 'NO'
$arr3[12,0] = '4567'
$arr3[12,1] = 4
$arr3[12,2] = ''
$arr3[12,3] = ''
$arr3[12,4] = ''
$arr3[12,5] = 'NO'
$arr3[13,0] = '4586'
$arr3[13,1] = ""
$arr3[13,2] = ''
$arr3[13,3] = ''
$arr3[13,4] = ''
$arr3[13,5] = 'NO'
$arr3[14,0] = '4590'
$arr3[14,1] = ""
$arr3[14,2] = ''
$arr3[14,3] = ''
$arr3[14,4] = ''
$arr3[14,5] = 'NO'
$arr3[15,0] = '4592'
$arr3[15,1] = 4
$arr3[15,2] = '0'
$arr3[15,3] = '0'
$arr3[15,4] = '4.09%'
$arr3[15,5] = 'NO'
$arr3[16,0] = '4634'
$arr3[16,1] = ""
$arr3[16,2] = ''
$arr3[16,3] = ''
$arr3[16,4] = ''
$arr3[16,5] = 'NO'
$arr3[17,0] = '4638'
$arr3[17,1] = 4
$arr3[17,2] = '6'
$arr3[17,3] = '1'
$arr3[17,4] = '36.13%'
$arr3[17,5] = 'NO'
$arr3[18,0] = '4686'
$arr3[18,1] = ""
$arr3[18,2] = ''
$arr3[18,3] = ''
$arr3[18,4] = ''
$arr3[18,5] = 'NO'
$arr3[19,0] = '4688'
$arr3[19,1] = 4
$arr3[19,2] = '2'
$arr3[19,3] = '0'
$arr3[19,4] = '15.38%'
$arr3[19,5] = 'NO'
$arr3[20,0] = '4690'
$arr3[20,1] = ""
$arr3[20,2] = ''
$arr3[20,3] = ''
$arr3[20,4] = ''
$arr3[20,5] = 'NO'
$battingExtra.Range("A1:F21").NumberFormat = "@"
$battingExtra.Range("B2:B21").NumberFormat = "General"
$battingExtra.Range("A1:F21").Value = $arr3

Write-Host "Edit complete."
